# Data update for past 3 days
# Appends rows 67-69 (2020-05-23, 2020-05-24, 2020-05-25) to the coronadata
# age/sex breakdown sheet, matching the existing table layout:
#   A=Date, B..U = age/sex bucket counts, V=M total, W=F total, X=%M, Y=%F

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 67; Date = 43974; Values = @(43, 69, 107, 261, 135, 234, 109, 189, 156, 225, 167, 182, 147, 177, 103, 146, 46, 121, 3, 7) },
    @{ Row = 68; Date = 43975; Values = @(43, 69, 107, 268, 136, 236, 109, 190, 157, 226, 168, 183, 147, 178, 105, 147, 47, 122, 3, 7) },
    @{ Row = 69; Date = 43976; Values = @(44, 69, 109, 272, 136, 240, 111, 192, 158, 227, 168, 185, 147, 180, 105, 152, 48, 123, 3, 7) }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    # Column A: date (inherits the existing yyyy-mm-dd number format from the row above)
    $ws.Cells.Item($r, 1).Value = $rowData.Date

    # Columns B..U: age/sex bucket counts
    $values = $rowData.Values
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $values[$i]
    }

    # Column V: sum of male buckets, W: sum of female buckets
    $ws.Range("V$r").Formula = "=SUM(B$r,D$r,F$r,H$r,J$r,L$r,N$r,P$r,R$r,T$r)"
    $ws.Range("W$r").Formula = "=SUM(C$r,E$r,G$r,I$r,K$r,M$r,O$r,Q$r,S$r,U$r)"

    # Columns X/Y: male/female percentage of the day's total
    $ws.Range("X$r").Formula = "=(V$r/(V$r+W$r))*100"
    $ws.Range("Y$r").Formula = "=(W$r/(V$r+W$r))*100"
}

# Match the author's final selection on the sheet
$ws.Range("V79").Select()
